$d = $word.ActiveDocument

# 1. Purchase Order Number validation message
$d.Content.Find.Execute(
    "Purchase Order Number must be alphanumeric only, use of special characters is prohibited.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Please enter Purchase Order Number containing letters, numbers, and hyphens only.",
    2)

# 2. Date Due validation message (mm/dd/yy -> mm/dd/yyyy)
$d.Content.Find.Execute(
    "Date Due must be of the format mm/dd/yy.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Date Due must be of the format mm/dd/yyyy.",
    2)

# 3. Date Received validation message (must of -> must be of, mm/dd/yy -> mm/dd/yyyy)
$d.Content.Find.Execute(
    "Date Received must of the format mm/dd/yy.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Date Received must be of the format mm/dd/yyyy.",
    2)
